$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45304) {
        $cell.Value2 = 45305
    }
}
